# Refresh the crypto price ticker data (Price / Volume(1h) columns, plus the
# BabyDogeCoin/RocketPoolETH row swap) to match the latest GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.023.30'
$ws.Range('E2').Value = '  +0.00%  '

$ws.Range('D3').Value = '1.832.67'
$ws.Range('E3').Value = '  +0.09%  '

$ws.Range('D4').Value = '''0.9978'
$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').Value = '''242.43'
$ws.Range('E5').Value = '  +0.43%  '

$ws.Range('D6').Value = '''0.6271'
$ws.Range('E6').Value = '  -4.12%  '

$ws.Range('D7').Value = '''0.9999'
$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('E8').Value = '  +3.54%  '

$ws.Range('D9').Value = '''0.2923'
$ws.Range('E9').Value = '  -0.43%  '

$ws.Range('D10').Value = '''22.57'
$ws.Range('E10').Value = '  -1.75%  '

$ws.Range('D11').Value = '''0.07718'
$ws.Range('E11').Value = '  +0.54%  '

$ws.Range('D12').Value = '1.830.63'
$ws.Range('E12').Value = '  -0.13%  '

$ws.Range('D13').Value = '''4.951'
$ws.Range('E13').Value = '  -0.69%  '

$ws.Range('D14').Value = '''0.6639'
$ws.Range('E14').Value = '  -0.36%  '

$ws.Range('D15').Value = '''0.00001023'
$ws.Range('E15').Value = '  +18.44%  '

$ws.Range('D16').Value = '''82.77'
$ws.Range('E16').Value = '  +0.59%  '

$ws.Range('D17').Value = '''6.049'
$ws.Range('E17').Value = '  -0.16%  '

$ws.Range('D18').Value = '29.030.38'
$ws.Range('E18').Value = '  +0.25%  '

$ws.Range('D19').Value = '''226.66'
$ws.Range('E19').Value = '  +1.05%  '

$ws.Range('D20').Value = '''12.34'
$ws.Range('E20').Value = '  -0.67%  '

$ws.Range('D21').Value = '''0.9989'
$ws.Range('E21').Value = '  -0.06%  '

$ws.Range('D22').Value = '''7.185'
$ws.Range('E22').Value = '  +1.07%  '

$ws.Range('D23').Value = '''0.9991'
$ws.Range('E23').Value = '  -0.04%  '

$ws.Range('D24').Value = '''158.57'
$ws.Range('E24').Value = '  +0.50%  '

$ws.Range('D25').Value = '''8.493'
$ws.Range('E25').Value = '  -0.14%  '

$ws.Range('D26').Value = '''0.1373'
$ws.Range('E26').Value = '  -0.58%  '

$ws.Range('D27').Value = '''17.88'
$ws.Range('E27').Value = '  -0.33%  '

$ws.Range('D28').Value = '''1.491'
$ws.Range('E28').Value = '  -0.95%  '

$ws.Range('D29').Value = '''4.094'
$ws.Range('E29').Value = '  -0.37%  '

$ws.Range('D30').Value = '''4.013'
$ws.Range('E30').Value = '  +0.09%  '

$ws.Range('D31').Value = '''1.188'
$ws.Range('E31').Value = '  -1.45%  '

$ws.Range('D32').Value = '''0.05237'
$ws.Range('E32').Value = '  -2.03%  '

$ws.Range('D33').Value = '''1.842'
$ws.Range('E33').Value = '  +0.62%  '

$ws.Range('D34').Value = '''0.7342'
$ws.Range('E34').Value = '  -1.33%  '

$ws.Range('D35').Value = '''1.139'
$ws.Range('E35').Value = '  -1.20%  '

$ws.Range('D36').Value = '''2.692'
$ws.Range('E36').Value = '  +1.93%  '

$ws.Range('D37').Value = '1.233.54'
$ws.Range('E37').Value = '  -4.77%  '

$ws.Range('D38').Value = '''2.755'
$ws.Range('E38').Value = '  +0.46%  '

$ws.Range('D39').Value = '''0.01785'
$ws.Range('E39').Value = '  -0.04%  '

$ws.Range('D40').Value = '''6.352'
$ws.Range('E40').Value = '  +0.12%  '

$ws.Range('D41').Value = '''0.8980'
$ws.Range('E41').Value = '  -0.01%  '

$ws.Range('D42').Value = '''0.9999'
$ws.Range('E42').Value = '  +0.10%  '

$ws.Range('D43').Value = '''102.03'
$ws.Range('E43').Value = '  -1.16%  '

$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = '1.981.20'
$ws.Range('E44').Value = '  -0.08%  '

$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = '''0.00000000123'
$ws.Range('E45').Value = '  +2.64%  '

$ws.Range('D46').Value = '''64.13'
$ws.Range('E46').Value = '  +0.17%  '

$ws.Range('D47').Value = '''0.5104'
$ws.Range('E47').Value = '  -0.70%  '

$ws.Range('D48').Value = '''0.4037'
$ws.Range('E48').Value = '  +1.31%  '

$ws.Range('D49').Value = '''8.895'
$ws.Range('E49').Value = '  +2.33%  '

$ws.Range('D50').Value = '''0.05737'
$ws.Range('E50').Value = '  -1.76%  '

$ws.Range('D51').Value = '''6.679'
$ws.Range('E51').Value = '  -0.31%  '
